$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume figures ("Updated symbol list" GitHub Actions commit).
# Cells D/E on each row store Price / Volume(1h) as plain text (inline strings),
# so NumberFormat is forced to Text ("@") before assignment to stop Excel from
# re-interpreting the numeric-looking / percent-looking strings as real numbers
# (which would silently drop trailing zeros / rescale percentages).

$updates = @{
    "D2" = "301.71"
    "D3" = "35.13"
    "E3" = "-0.47%"
    "D4" = "5.041"
    "E4" = "-1.84%"
    "D5" = "0.07982"
    "D6" = "1.916"
    "E6" = "-10.19%"
    "D7" = "7.804"
    "E7" = "-2.39%"
    "E8" = "8.97%"
    "D9" = "0.9229"
    "E9" = "-0.71%"
    "D10" = "0.1323"
    "E10" = "31.74%"
    "D11" = "0.1847"
    "E11" = "-1.50%"
    "D12" = "0.09641"
    "E12" = "6.19%"
    "D13" = "0.03587"
    "E13" = "-0.35%"
    "D14" = "0.09857"
    "E14" = "-0.33%"
    "D15" = "0.001387"
    "E15" = "-3.19%"
    "D16" = "0.005814"
    "E16" = "1.61%"
    "D17" = "3.504"
    "E17" = "0.79%"
    "D18" = "4.046"
    "E18" = "-2.44%"
    "D19" = "0.3401"
    "E19" = "-0.20%"
    "D20" = "0.1301"
    "E20" = "-2.23%"
    "D21" = "5.060"
    "E21" = "-0.94%"
    "D22" = "0.2400"
    "E22" = "8.09%"
    "D23" = "0.04500"
    "E23" = "-1.46%"
    "D24" = "0.001216"
    "E24" = "-2.48%"
    "D25" = "0.004785"
    "E25" = "1.67%"
    "E26" = "-0.14%"
    "D27" = "0.0003002"
    "E27" = "-33.36%"
    "D39" = "0.01882"
    "E39" = "-3.51%"
    "D40" = "0.04715"
    "E40" = "-2.63%"
    "D41" = "0.007492"
    "E41" = "-3.19%"
    "D42" = "0.01021"
    "E42" = "30.15%"
    "D43" = "0.1325"
    "E43" = "-4.67%"
    "E44" = "-0.11%"
    "D45" = "0.01062"
    "E45" = "-8.65%"
    "D46" = "0.00006233"
    "E46" = "-5.69%"
    "D47" = "0.00000000750"
    "E47" = "-0.04%"
    "E48" = "72.45%"
    "E49" = "-12.47%"
    "D50" = "0.00002101"
    "E50" = "-0.04%"
    "D51" = "0.0002001"
    "E51" = "-0.04%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

